$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new hashcode value (column B)
$updates = @{
    44  = "a47a153b559b63e964da9aed9585b7a1"
    89  = "22a3d0d0ffbd6966636b8cc771a4a7f3"
    99  = "2301c9604a28fed3e5d97db82d9b9f3e"
    110 = "6850bac7505a0d9eeab4940d7358fc07"
    154 = "57735bff7c761462623e6fe87e0de54c"
    160 = "359d356b5210f737e5dbecd965fc1b23"
    168 = "a3f28045d5ed069e7720a844e2fed592"
    222 = "fcf9d67bfccbed63010e25156f3e1f83"
    229 = "415c536bd2b2a1ac6ddb0bcee24d48c4"
    278 = "66291cd0fe7fe3f41ad6cc951fb55cb7"
    330 = "bd173cda6e98b8b319511b4ab70dbbeb"
    335 = "b9d82cf15c770f67d2ab4afad88a710e"
    411 = "3a43ae4697bc716a57c569af214f226e"
    444 = "da36aaf095228edc5d62cbe6d1c6a17a"
    448 = "13ad8e9c3d9cf959fdbaf7938c139475"
    523 = "700ef0f36fdad29e5cb174cccc20b89d"
    561 = "52e2c43fa48846733196fed4fa7f88dd"
    574 = "9db01d3d4a4b1b90e07d70b68338ffb5"
    592 = "6b6213daff9d5f4fe046f1a5a5c5c70c"
    764 = "450c7df6e48c330d652f26b07d339735"
    769 = "3ef00d0b56e28622e10a6591a238e6bc"
    776 = "9dcd3d5ac9464a0a67843c938917f803"
    794 = "551495b228c1ec24f7f68f3754f9c716"
    833 = "080f53956c653b496a0e2afea805c3dd"
    835 = "c3f7248a6171c0e2ee7f819a9bebd5fe"
    863 = "49a0c6b32d97d34e893f73194ae95d1b"
    877 = "174be31a5a3b3ebd01a92b82e698cf8c"
    882 = "d878f735a89572d2273c1e98708e28dd"
    913 = "9eeee1a400bfb91f95e3f221321ea772"
    937 = "53b155dc7bbaba9d7d889111a44dfaf0"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
